# Styling, help, log, dll
#
# - Rename "Sheet1" to "Sheet 1"
# - Make "Sheet 1" the active sheet/tab (was STRESS) and set its selection to G20
# - As a consequence the previously active "TEST" tab is no longer tabSelected

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Name = "Sheet 1"

$sheet1.Activate()
[void]$sheet1.Range("G20").Select()
